# Change the workbook's default ("Normal") font from Calibri to Arial.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Styles("Normal").Font.Name = "Arial"

# Enter the text "bbb" into the first populated cell (C1) and select it,
# matching the new active cell / dimension of the sheet.
$ws.Range("C1").Value = "bbb"
$ws.Range("C1").Select()

$excel.ActiveWindow.WindowState = -4143
